$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Fill in newly-documented "LCD Display State" cells (column D) ---

# D29: handling indications state
$ws.Range("D29").Value = "DISPLAY_ROW_CONNECTION - handling indications"

# D13: connected state (copy formatting from a sibling cell that already
# carries style 29, then set the value)
$ws.Range("E7").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").Value = "DISPLAY_ROW_CONNECTION - connected, DISPLAY_ROW_BTADDR2 - server address"

# D36: temp value displayed
$ws.Range("D36").Value = "DISPLAY_ROW_TEMPVALUE - temp value"

# D7: discovering state (also needs style 29, matching E7's formatting)
$ws.Range("E7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").Value = "DISPLAY_ROW_CONNECTION - discovering, DISPLAY_ROW_NAME - client, DISPLAY_ROW_BTADDR - client address, DISPLAY_ROW_ASSIGNMENT - A7"

$excel.CutCopyMode = 0

# --- Update the sheet view (scroll position, zoom, selection) ---
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 116
$ws.Range("E17").Select() | Out-Null
